$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.361.21"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "1.655.92"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'213.43"
$ws.Range("E5").Value = "  -0.67%  "
$ws.Range("D6").Value = "'0.515"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("D8").Value = "'23.62"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "1.890.31"
$ws.Range("E12").Value = "  -0.75%  "
$ws.Range("D13").Value = "1.690.71"
$ws.Range("E13").Value = "  +1.32%  "
$ws.Range("E14").Value = "  +3.72%  "
$ws.Range("E15").Value = "  -1.77%  "
$ws.Range("D16").Value = "'65.78"
$ws.Range("E16").Value = "  -0.49%  "
$ws.Range("D17").Value = "27.361.30"
$ws.Range("E17").Value = "  -1.70%  "
$ws.Range("D18").Value = "'231.94"
$ws.Range("E18").Value = "  -7.80%  "
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  -1.24%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  -2.45%  "
$ws.Range("D23").Value = "'9.36"
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  -1.46%  "
$ws.Range("D25").Value = "'146.94"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  -2.66%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.19%  "
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("E32").Value = "  -1.49%  "
$ws.Range("D33").Value = "1.448.14"
$ws.Range("E33").Value = "  +1.23%  "
$ws.Range("D34").Value = "'3.14"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +0.34%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("D37").Value = "'0.909"
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("E38").Value = "  -2.17%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D42").Value = "'5.55"
$ws.Range("E42").Value = "  +2.88%  "
$ws.Range("E43").Value = "  -6.71%  "
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'0.788"
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.798.02"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("D47").Value = "'1.71"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("D48").Value = "'88.10"
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("D49").Value = "0.0₆0107"
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "'7.75"
$ws.Range("E51").Value = "  -0.86%  "
